$d = $word.ActiveDocument

# Helper: turn the (empty) paragraph currently at 1-based index $idx into a
# paragraph that carries only paragraph-mark run formatting (sz=28, szCs=28,
# lang=nb-NO) in its <w:pPr><w:rPr>, with no literal run in the body.
#
# We briefly type a placeholder character into the paragraph so the engine
# has real run content to apply character formatting to (Font.Size /
# Font.SizeBi / LanguageID are no-ops on a completely empty, "just the
# paragraph mark" range), then delete that character again. Deleting it
# collapses the run but keeps the formatting recorded on the paragraph mark
# (<w:pPr><w:rPr>), matching how Word stores an empty, pre-formatted
# paragraph.
function Format-TrailingParagraph($idx) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertBefore("X")

    $p = $d.Paragraphs.Item($idx)
    $p.Range.Font.Size = 14
    $p.Range.Font.SizeBi = 14
    $p.Range.LanguageID = "nb-NO"

    $charRange = $d.Range($p.Range.Start, $p.Range.Start + 1)
    $charRange.Delete()

    $p = $d.Paragraphs.Item($idx)
    $p.Range.LanguageID = "nb-NO"
}

# Two new empty paragraphs are added right after the "Løkker..." paragraph,
# i.e. right before the trailing empty paragraph that precedes the section
# break. Insert them one at a time, immediately before whatever paragraph is
# currently last, so each new paragraph lands in the right spot.
for ($n = 0; $n -lt 2; $n++) {
    $lastIdx = $d.Paragraphs.Count
    $lastPara = $d.Paragraphs.Item($lastIdx)
    $lastPara.Range.InsertParagraphBefore()

    $newIdx = $d.Paragraphs.Count - 1
    Format-TrailingParagraph($newIdx)
}
